# Update "想去人数" (column F) counts across sheets, per commit
# "Update gh-pages to output generated at 456a3b4".
# Sheet order in workbook: 1=展览, 2=演出, 3=本地生活 (unaffected), 4=全部类型

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 15097
$ws.Range("F3").Value = 19288
$ws.Range("F15").Value = 236
$ws.Range("F16").Value = 70
$ws.Range("F17").Value = 1491
$ws.Range("F20").Value = 103
$ws.Range("F22").Value = 8078
$ws.Range("F25").Value = 6
$ws.Range("F27").Value = 1257
$ws.Range("F28").Value = 6
$ws.Range("F30").Value = 6096
$ws.Range("F31").Value = 123
$ws.Range("F36").Value = 5504
$ws.Range("F37").Value = 1008
$ws.Range("F38").Value = 21
$ws.Range("F40").Value = 55

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 23

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 15097
$ws.Range("F3").Value = 19288
$ws.Range("F15").Value = 236
$ws.Range("F16").Value = 70
$ws.Range("F17").Value = 1491
$ws.Range("F21").Value = 103
$ws.Range("F23").Value = 8078
$ws.Range("F26").Value = 6
$ws.Range("F28").Value = 1257
$ws.Range("F29").Value = 6
$ws.Range("F31").Value = 23
$ws.Range("F33").Value = 6096
$ws.Range("F34").Value = 123
$ws.Range("F39").Value = 5504
$ws.Range("F40").Value = 1008
$ws.Range("F41").Value = 21
$ws.Range("F43").Value = 55
